# Updates the crypto price table (columns D-G, rows 2-51) to the
# Feb 5 2023 00:30 UTC snapshot values pulled by the GitHub Actions job.
#
# Cells in this sheet are stored as text (inline strings) even though many
# of the values look numeric (prices, percentages, the "0" hour, etc.).
# Setting `.Value` directly on a numeric-looking string makes Excel coerce
# it into a real number (and reformat percentages), which would not match
# the original text-cell layout. Forcing the cell to the "@" (Text) number
# format before the assignment keeps the literal text, and then clearing
# the format afterwards (`ClearFormats`) drops the now-unneeded "Text"
# style so the cell's style index is left exactly as it was before (no
# explicit `s=` attribute), matching how the rest of the sheet looks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Value) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.ClearFormats()
}

$rowData = @(
    @{ Row = 2; D = "330.32"; E = "-0.81%"; F = "5-2-2023"; G = "0" }
    @{ Row = 3; D = "41.11"; E = "-0.02%"; F = "5-2-2023"; G = "0" }
    @{ Row = 4; D = "5.668"; E = "-1.43%"; F = "5-2-2023"; G = "0" }
    @{ Row = 5; D = "0.08366"; E = "1.93%"; F = "5-2-2023"; G = "0" }
    @{ Row = 6; D = "8.812"; E = "0.78%"; F = "5-2-2023"; G = "0" }
    @{ Row = 7; D = "4.510"; E = "-0.24%"; F = "5-2-2023"; G = "0" }
    @{ Row = 8; D = "1.987"; E = "-2.84%"; F = "5-2-2023"; G = "0" }
    @{ Row = 9; D = "2.924"; E = "-2.49%"; F = "5-2-2023"; G = "0" }
    @{ Row = 10; D = "0.9254"; E = "0.45%"; F = "5-2-2023"; G = "0" }
    @{ Row = 11; D = "0.1248"; E = "0.49%"; F = "5-2-2023"; G = "0" }
    @{ Row = 12; D = "0.1962"; E = "0.57%"; F = "5-2-2023"; G = "0" }
    @{ Row = 13; D = "0.09427"; E = "0.24%"; F = "5-2-2023"; G = "0" }
    @{ Row = 14; D = "0.03956"; E = "8.06%"; F = "5-2-2023"; G = "0" }
    @{ Row = 15; D = "0.1065"; E = "0.92%"; F = "5-2-2023"; G = "0" }
    @{ Row = 16; D = "0.001319"; E = "1.71%"; F = "5-2-2023"; G = "0" }
    @{ Row = 17; D = "0.006108"; E = "-1.62%"; F = "5-2-2023"; G = "0" }
    @{ Row = 18; D = "3.437"; E = "1.50%"; F = "5-2-2023"; G = "0" }
    @{ Row = 19; D = $null; E = $null; F = "5-2-2023"; G = "0" }
    @{ Row = 20; D = "9.053"; E = "9.26%"; F = "5-2-2023"; G = "0" }
    @{ Row = 21; D = "0.1375"; E = "-3.07%"; F = "5-2-2023"; G = "0" }
    @{ Row = 22; D = "0.2639"; E = "-0.56%"; F = "5-2-2023"; G = "0" }
    @{ Row = 23; D = "0.04438"; E = "-0.07%"; F = "5-2-2023"; G = "0" }
    @{ Row = 24; D = "0.001251"; E = "-0.90%"; F = "5-2-2023"; G = "0" }
    @{ Row = 25; D = "0.004382"; E = "1.87%"; F = "5-2-2023"; G = "0" }
    @{ Row = 26; D = "0.0001196"; E = "-3.60%"; F = "5-2-2023"; G = "0" }
    @{ Row = 27; D = "0.0004010"; E = "0.43%"; F = "5-2-2023"; G = "0" }
    @{ Row = 28; D = $null; E = $null; F = "5-2-2023"; G = "0" }
    @{ Row = 29; D = $null; E = $null; F = "5-2-2023"; G = "0" }
    @{ Row = 30; D = $null; E = $null; F = "5-2-2023"; G = "0" }
    @{ Row = 31; D = $null; E = $null; F = "5-2-2023"; G = "0" }
    @{ Row = 32; D = $null; E = $null; F = "5-2-2023"; G = "0" }
    @{ Row = 33; D = $null; E = $null; F = "5-2-2023"; G = "0" }
    @{ Row = 34; D = $null; E = $null; F = "5-2-2023"; G = "0" }
    @{ Row = 35; D = $null; E = $null; F = "5-2-2023"; G = "0" }
    @{ Row = 36; D = $null; E = $null; F = "5-2-2023"; G = "0" }
    @{ Row = 37; D = $null; E = $null; F = "5-2-2023"; G = "0" }
    @{ Row = 38; D = $null; E = $null; F = "5-2-2023"; G = "0" }
    @{ Row = 39; D = "0.02808"; E = "0.08%"; F = "5-2-2023"; G = "0" }
    @{ Row = 40; D = "0.05510"; E = "0.12%"; F = "5-2-2023"; G = "0" }
    @{ Row = 41; D = "0.007935"; E = "3.53%"; F = "5-2-2023"; G = "0" }
    @{ Row = 42; D = "0.1430"; E = "0.54%"; F = "5-2-2023"; G = "0" }
    @{ Row = 43; D = "0.009097"; E = "-8.69%"; F = "5-2-2023"; G = "0" }
    @{ Row = 44; D = "0.002179"; E = "2.26%"; F = "5-2-2023"; G = "0" }
    @{ Row = 45; D = $null; E = "-14.71%"; F = "5-2-2023"; G = "0" }
    @{ Row = 46; D = "0.00007205"; E = "6.69%"; F = "5-2-2023"; G = "0" }
    @{ Row = 47; D = "0.00000000754"; E = "0.40%"; F = "5-2-2023"; G = "0" }
    @{ Row = 48; D = "0.003553"; E = "17.29%"; F = "5-2-2023"; G = "0" }
    @{ Row = 49; D = "0.002290"; E = "0.33%"; F = "5-2-2023"; G = "0" }
    @{ Row = 50; D = "0.00002111"; E = "0.40%"; F = "5-2-2023"; G = "0" }
    @{ Row = 51; D = "0.0002011"; E = "0.40%"; F = "5-2-2023"; G = "0" }
)

foreach ($entry in $rowData) {
    $r = $entry.Row
    if ($null -ne $entry.D) {
        Set-TextValue $ws.Cells.Item($r, 4) $entry.D
    }
    if ($null -ne $entry.E) {
        Set-TextValue $ws.Cells.Item($r, 5) $entry.E
    }
    Set-TextValue $ws.Cells.Item($r, 6) $entry.F
    Set-TextValue $ws.Cells.Item($r, 7) $entry.G
}
